# Arbeitsaufzeichnung von Aron fertig
# Appends the remaining work-log entries (rows 35-50), widens column A,
# extends the duration total formula, and adds the small left-border
# accent on the "Array der Gesichtspunkte..." row - matching the
# author's final commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-look-alike string ("DD.MM.YYYY" with both parts
# <=12) as literal text. Excel's Value setter auto-converts those into
# date serials (ambiguous day/month), so the cell is briefly forced to
# Text format for the assignment and then returned to the Normal style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- New data rows (written in the exact order the cell text was first
# authored, so new shared-string entries land in the same sequence as
# the target workbook) -------------------------------------------------

# Row 36 header/body text entered before row 35's in the source history.
$ws.Range("A36").Value = "Fehler in Log speichern"
$ws.Range("A35").Value = "Status-Präsentation vorbereiten"
$ws.Range("B35").Value = "13.12.2019-16.12.2019"
$ws.Range("B36").Value = "23.12.2019-26.12.2019"
$ws.Range("C35").Value = 4
$ws.Range("C36").Value = 5

$ws.Range("A37").Value = "Array der Gesichtspunkte in zwei aufteilen. Das erste ist nur für die x-Achse und das zweite für die y-Achse"
$ws.Range("B37").Value = "04.01.2020-07.01.2020"
$ws.Range("C37").Value = 7

$ws.Range("A38").Value = "1 Betreuertreffung"
$ws.Range("B38").Value = "14.01.2020"
$ws.Range("C38").Value = 1

$ws.Range("A39").Value = "Dokumentation Korrektur"
$ws.Range("B39").Value = "17.01.2020"
$ws.Range("C39").Value = 4

$ws.Range("A40").Value = "Betreuertreffung"
$ws.Range("B40").Value = "25.01.2020"
$ws.Range("C40").Value = 1

$ws.Range("A41").Value = "3 Prototyps gemacht für Backups, wenn ein nicht funktioniert"
$ws.Range("B41").Value = "25.01.2020-27-01.2020"
$ws.Range("C41").Value = 6

$ws.Range("A42").Value = "Dokumentation verarbeiten, weiter mit der Korrektur, Struktur der Dokumentation verbessern"
$ws.Range("B42").Value = "31.01.2020"
$ws.Range("C42").Value = 4

$ws.Range("A43").Value = "In Registrierung-Teil arbeiten, verschiedene Skripts miteinander verbinden, die richtigen Konsole-Parameter von Skript zu Skript übergeben"
$ws.Range("B43").Value = "01.02.2020-02.02.2020"
$ws.Range("C43").Value = 4

$ws.Range("A44").Value = "Automatische Konfiguration von Opencv in Raspberry PI"
Set-TextValue $ws.Range("B44") "03.02.2020"
$ws.Range("C44").Value = 3

$ws.Range("A45").Value = "Log (in Datenbank und in einer Text-File)"
Set-TextValue $ws.Range("B45") "05.02.2020"
$ws.Range("C45").Value = 5

$ws.Range("A46").Value = "Anpassung des Registrierungsteils, weil in der DB etwas geändert hat (von zwei Tabellen, jetzt nur eine)"
$ws.Range("B46").Value = "07.02.2020,08.02.2020"
$ws.Range("C46").Value = 3

$ws.Range("A47").Value = "Das Problem mit dem Zugriff auf Elementen des numpy-Arrays analysieren und lösen"
Set-TextValue $ws.Range("B47") "10.02.2020"
$ws.Range("C47").Value = 3

# Row 48: the date (B) was authored before the task text (A).
$ws.Range("B48").Value = "11.02.2020-15.02.2020"
$ws.Range("A48").Value = "Testen des Registrierungsteil"
$ws.Range("C48").Value = 10

$ws.Range("A49").Value = "In der Finalabgabe der Dokumentation arbeiten"
$ws.Range("B49").Value = "16.02.2020-20.02.2020, 30.02.2020-03.03.2020"
$ws.Range("C49").Value = 15

$ws.Range("A50").Value = "Die End-Präsentation von der Diplomarbeit "
$ws.Range("B50").Value = "04.03.2020-06.03.2020"
$ws.Range("C50").Value = 5

# --- Small left-border accent on row 37 (matches the thin left border
# the author added around the "Array der Gesichtspunkte..." entry) ----
$ws.Range("A37").Borders.Item(7).LineStyle = 1
$ws.Rows.Item(37).RowHeight = 13

# --- Widen column A to fit the longer task descriptions ---------------
$ws.Columns.Item(1).ColumnWidth = 123

# --- Extend the duration total to cover the new rows -------------------
$ws.Range("E3").Formula = "=SUM(C3:C50)"

# --- Keep the header selection where the author left it ----------------
$ws.Range("C5").Select()
